$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C2:C8").NumberFormat = "[$-1010000]yyyy/mm/dd;@"
$ws.Range("C2").Value = 36620
$ws.Range("C3").Value = 36620
$ws.Range("C4").Value = 36621
$ws.Range("C5").Value = 36622
$ws.Range("C6").Value = 36623
$ws.Range("C7").Value = 36624
$ws.Range("C8").Value = 36625

$ws.Range("M18").Select() | Out-Null
